$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Overwrite the test result for the existing user (Tanja) in C3
$ws.Range("C3").Value = 8

# Clear out the duplicate/stale rows 5-7 (contents only, rows remain)
$ws.Range("A5:C7").ClearContents()

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("E14").Select()
